$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 79
$ws1.Range("F4").Value = 1515
$ws1.Range("F5").Value = 583
$ws1.Range("F6").Value = 1073
$ws1.Range("F7").Value = 11123
$ws1.Range("F10").Value = 322
$ws1.Range("F11").Value = 1069
$ws1.Range("F12").Value = 759
$ws1.Range("F13").Value = 12246
$ws1.Range("F14").Value = 12814
$ws1.Range("F16").Value = 126
$ws1.Range("F21").Value = 49

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 79
$ws4.Range("F5").Value = 1515
$ws4.Range("F6").Value = 583
$ws4.Range("F7").Value = 1073
$ws4.Range("F8").Value = 11123
$ws4.Range("F11").Value = 322
$ws4.Range("F12").Value = 1069
$ws4.Range("F13").Value = 759
$ws4.Range("F14").Value = 12246
$ws4.Range("F15").Value = 12814
$ws4.Range("F17").Value = 126
$ws4.Range("F22").Value = 49
